# Schema Migration Template - update teradata host IP
#
# The source "host" value (key "host" in row 10 of the
# "Schema migration input template" sheet) is stale; refresh it to the
# new teradata host IP and make sure the cell reads cleanly (wrap +
# vertically centered) now that the value is longer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schema migration input template")
$ws.Activate()

$hostCell = $ws.Range("B10")
$hostCell.Select()

# Drop any stale direct formatting on the cell before re-applying the
# look we want, then write the new host IP.
$hostCell.ClearFormats()
$hostCell.Value = "34.31.79.171"

# Keep the value readable: wrap text and center it vertically.
$hostCell.WrapText = $true
$hostCell.VerticalAlignment = -4108
